# Improve error message display / Update dependencies to their latest versions.
# Append one new data row to each of the four worksheets, matching the
# existing layout (time / length / ID / actual-length / checksum + their
# decimal counterparts).

$wb = $excel.ActiveWorkbook

function Add-DataRow {
    param($ws, $row, $timeValue, $colB, $colC, $colD, $colE, $colF, $colG, $colH, $colI)

    # Column A: date/time serial, formatted like the rest of the column.
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 1).Value = $timeValue

    # Columns B-E: hex-byte strings stored as text.
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $colD
    $ws.Cells.Item($row, 5).Value = $colE

    # Column F: plain numeric.
    $ws.Cells.Item($row, 6).Value = $colF

    # Column G: numeric decimal id -- some values overflow double precision
    # and must be kept as literal text (matches existing rows in the sheet).
    if ($colG -is [string]) {
        $ws.Cells.Item($row, 7).NumberFormat = "@"
        $ws.Cells.Item($row, 7).Value = $colG
        $ws.Cells.Item($row, 7).Style = "Normal"
    } else {
        $ws.Cells.Item($row, 7).Value = $colG
    }

    # Columns H-I: plain numeric.
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
}

$gNum = [double]"5.68631262647114e+23"

# Sheet "ROW50-FE-LIFTER": new row 42
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
Add-DataRow `
    $ws1 `
    42 `
    45745.17660361111 `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," `
    "0x01,0x6a" `
    "0xe" `
    400 `
    $gNum `
    362 `
    14

# Sheet "ROW50-MID-LIFTER": new row 44
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
Add-DataRow `
    $ws2 `
    44 `
    45745.1428125 `
    "0x01,0x90 " `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," `
    "0x01,0x72" `
    "0x19" `
    400 `
    "568631262647113771663628" `
    370 `
    25

# Sheet "ROW11-FE-LIFTER": new row 42
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
Add-DataRow `
    $ws3 `
    42 `
    45745.19591043981 `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," `
    "0x01,0x6a" `
    "0x14" `
    400 `
    $gNum `
    362 `
    20

# Sheet "ROW11-MID-LIFTER": new row 42
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
Add-DataRow `
    $ws4 `
    42 `
    45745.33860818287 `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," `
    "0x01,0x72" `
    "0x19" `
    400 `
    $gNum `
    370 `
    25
